# Junction_Flooding_92.xlsx edit:
#  1. Row 5 values are re-expressed at "custom accuracy" (rounded to 2 decimals).
#  2. Row 6 (the last data row) is removed entirely -> dimension shrinks A1:AH6 -> A1:AH5.
#  3. A handful of column widths are trimmed by 1 character (8 -> 7, or 9 -> 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Row 5: rewrite values rounded to 2 decimal places ---------------
$ws.Range("B5").Value  = 14.6
$ws.Range("C5").Value  = 10.67
$ws.Range("D5").Value  = 0.68
$ws.Range("E5").Value  = 30.98
$ws.Range("F5").Value  = 26.09
$ws.Range("G5").Value  = 11.33
$ws.Range("H5").Value  = 41.4
$ws.Range("I5").Value  = 17.23
$ws.Range("J5").Value  = 7.74
$ws.Range("K5").Value  = 11.45
$ws.Range("L5").Value  = 12.58
$ws.Range("M5").Value  = 13.5
$ws.Range("N5").Value  = 3.6
$ws.Range("O5").Value  = 10.88
$ws.Range("P5").Value  = 16.46
$ws.Range("Q5").Value  = 9.01
$ws.Range("R5").Value  = 0.36
$ws.Range("S5").Value  = 0.42
$ws.Range("T5").Value  = 164.15
$ws.Range("U5").Value  = 31.43
$ws.Range("V5").Value  = 10.54
$ws.Range("W5").Value  = 21.14
$ws.Range("X5").Value  = 11.46
$ws.Range("Y5").Value  = 1.45
$ws.Range("Z5").Value  = 20.31
$ws.Range("AA5").Value = 9.16
$ws.Range("AB5").Value = 8.2
$ws.Range("AC5").Value = 10.03
$ws.Range("AD5").Value = 13.23
$ws.Range("AE5").Value = 0.57
$ws.Range("AF5").Value = 37.49
$ws.Range("AG5").Value = 5.79
$ws.Range("AH5").Value = 12.84

# --- 2. Delete row 6 entirely (last row of data) -------------------------
$ws.Rows.Item(6).Delete()

# --- 3. Column width tweaks -----------------------------------------------
# Excel's ColumnWidth (object model, "characters") differs from the raw
# OOXML <col width="..."> value by the fixed 5px/MDW padding constant
# (5/6 of a character for the default Calibri 11 font used here), so the
# COM-level width must be requested as (target OOXML width - 5/6) to land
# on the exact integer width recorded in the sheet XML.
$pad = 5.0 / 6.0

$ws.Columns.Item(2).ColumnWidth  = 7 - $pad   # B: 8 -> 7
$ws.Columns.Item(3).ColumnWidth  = 7 - $pad   # C: 8 -> 7
$ws.Columns.Item(11).ColumnWidth = 7 - $pad   # K: 8 -> 7
$ws.Columns.Item(12).ColumnWidth = 7 - $pad   # L: 8 -> 7
$ws.Columns.Item(13).ColumnWidth = 7 - $pad   # M: 8 -> 7
$ws.Columns.Item(15).ColumnWidth = 7 - $pad   # O: 8 -> 7
$ws.Columns.Item(16).ColumnWidth = 7 - $pad   # P: 8 -> 7
$ws.Columns.Item(20).ColumnWidth = 8 - $pad   # T: 9 -> 8
$ws.Columns.Item(22).ColumnWidth = 7 - $pad   # V: 8 -> 7
$ws.Columns.Item(24).ColumnWidth = 7 - $pad   # X: 8 -> 7
$ws.Columns.Item(26).ColumnWidth = 7 - $pad   # Z: 8 -> 7
$ws.Columns.Item(29).ColumnWidth = 7 - $pad   # AC: 8 -> 7
$ws.Columns.Item(30).ColumnWidth = 7 - $pad   # AD: 8 -> 7
$ws.Columns.Item(34).ColumnWidth = 7 - $pad   # AH: 8 -> 7
